# Rename the original sheet: "LayerCentrality" -> "Layer Centrality"
$wb = $excel.ActiveWorkbook
$base = $wb.Worksheets.Item(1)
$base.Name = "Layer Centrality"

# Last data row in the base sheet (row 1 = header, rows 2..62 = members, row 63 = overall mean)
$lastDataRow = 62

# Cluster definitions, in the order the tabs must appear.
# clusterId  -> per-cluster mean values for columns B..G (literal, matches source data exactly)
$clusterIds = @("0", "1", "2", "3", "-1")
$clusterMeans = @{
    "0"  = @(10.4425, 29.9, 14.18, 19.82916666666667, 25.65083333333333, 2.186666666666667)
    "1"  = @(0.1637931034482759, 20.51896551724138, 18.36, 29.70724137931035, 31.25137931034482, 1.848275862068966)
    "2"  = @(16.22, 3.7, 18.169, 29.897, 32.017, 2.061)
    "3"  = @(0.202, 5.465999999999999, 0.8560000000000001, 40.27999999999999, 53.2, 1.278)
    "-1" = @(9.720000000000001, 6.56, 0.5820000000000001, 43.50599999999999, 39.638, 1.266)
}

$prevSheet = $base
foreach ($cid in $clusterIds) {
    # Duplicate the fully-formatted base sheet (keeps styles, borders and the
    # conditional-formatting rules/dxfs intact) right after the previous tab.
    $base.Copy($null, $prevSheet)
    $new = $wb.Worksheets.Item($prevSheet.Index + 1)
    $new.Name = "cluster_$cid"

    # Overwrite the final (mean) row with this cluster's own averages.
    $means = $clusterMeans[$cid]
    $new.Cells.Item($lastDataRow + 1, 2).Value2 = $means[0]
    $new.Cells.Item($lastDataRow + 1, 3).Value2 = $means[1]
    $new.Cells.Item($lastDataRow + 1, 4).Value2 = $means[2]
    $new.Cells.Item($lastDataRow + 1, 5).Value2 = $means[3]
    $new.Cells.Item($lastDataRow + 1, 6).Value2 = $means[4]
    $new.Cells.Item($lastDataRow + 1, 7).Value2 = $means[5]
    $new.Cells.Item($lastDataRow + 1, 8).Value2 = [double]$cid

    # Drop every member row that doesn't belong to this cluster (bottom-up so
    # row numbers of not-yet-visited rows stay valid), counting the survivors
    # so the conditional-formatting range can be re-pointed afterwards.
    $kept = 0
    for ($r = $lastDataRow; $r -ge 2; $r--) {
        $h = $new.Cells.Item($r, 8).Value2
        if ($h -ne [double]$cid) {
            $new.Rows.Item($r).Delete()
        } else {
            $kept = $kept + 1
        }
    }

    # Re-point the conditional formatting range to the surviving member rows
    # (row 2 .. row (kept+1)); the mean row right below is excluded, just like
    # in the source sheet. This keeps reusing the same dxfs - no new styles
    # are created.
    $lastRow = $kept + 1
    $new.Range("B2:F2").FormatConditions.Item(1).ModifyAppliesToRange($new.Range("B2:F$lastRow"))

    $prevSheet = $new
}
